$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich-text (shared string) edits: Volume number and report week dates ---
$a8 = $ws.Range("A8")
$full = $a8.Value2
$idx = $full.IndexOf("51")
$a8.Characters($idx + 1, 2).Text = "52"

$c9 = $ws.Range("C9")
$full = $c9.Value2
$idx = $full.IndexOf("12/16/2024")
$c9.Characters($idx + 1, 10).Text = "12/23/2024"
$full = $c9.Value2
$idx = $full.IndexOf("12/22/2024")
$c9.Characters($idx + 1, 10).Text = "12/29/2024"

# --- Cell whose type flips from numeric to text (set text value first, then copy format from an already-text cell in the same column style family) ---
$ws.Range("C17").Value2 = "'0"
$ws.Range("A17").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null

# --- Cells whose type flips from text to numeric (copy format from an already-numeric cell in the same column style family, then set numeric value) ---
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value2 = 1
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value2 = 2
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value2 = -50

# --- Updated crime-data figures (numeric cells, value only) ---
# Row 16
$ws.Range("C16").Value2 = 3
$ws.Range("D16").Value2 = 6
$ws.Range("E16").Value2 = -50
$ws.Range("G16").Value2 = 12
$ws.Range("H16").Value2 = 16.666666666666
$ws.Range("I16").Value2 = 151
$ws.Range("J16").Value2 = 135
$ws.Range("K16").Value2 = 11.851851851851
$ws.Range("L16").Value2 = 12.686567164179
$ws.Range("M16").Value2 = 58.947368421052
$ws.Range("N16").Value2 = -82.543352601156

# Row 17
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = -100
$ws.Range("F17").Value2 = 11
$ws.Range("G17").Value2 = 17
$ws.Range("H17").Value2 = -35.294117647058
$ws.Range("J17").Value2 = 148
$ws.Range("K17").Value2 = 29.729729729729
$ws.Range("N17").Value2 = 3.783783783783

# Row 18
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 3
$ws.Range("F18").Value2 = 11
$ws.Range("G18").Value2 = 11
$ws.Range("H18").Value2 = 0
$ws.Range("I18").Value2 = 190
$ws.Range("J18").Value2 = 198
$ws.Range("K18").Value2 = -4.040404040404
$ws.Range("L18").Value2 = -26.640926640926
$ws.Range("M18").Value2 = 11.764705882352
$ws.Range("N18").Value2 = -79.074889867841

# Row 19
$ws.Range("C19").Value2 = 15
$ws.Range("D19").Value2 = 20
$ws.Range("E19").Value2 = -25
$ws.Range("F19").Value2 = 107
$ws.Range("G19").Value2 = 104
$ws.Range("H19").Value2 = 2.884615384615
$ws.Range("I19").Value2 = 1194
$ws.Range("J19").Value2 = 1315
$ws.Range("K19").Value2 = -9.201520912547
$ws.Range("L19").Value2 = -2.530612244897
$ws.Range("M19").Value2 = 12.112676056338
$ws.Range("N19").Value2 = -69.043297899922

# Row 20
$ws.Range("D20").Value2 = 1
$ws.Range("F20").Value2 = 2
$ws.Range("H20").Value2 = -75
$ws.Range("J20").Value2 = 78
$ws.Range("K20").Value2 = -43.589743589743
$ws.Range("N20").Value2 = -95.072788353863

# Row 21
$ws.Range("C21").Value2 = 21
$ws.Range("D21").Value2 = 31
$ws.Range("E21").Value2 = -32.258064516129
$ws.Range("F21").Value2 = 146
$ws.Range("G21").Value2 = 152
$ws.Range("H21").Value2 = -3.947368421052
$ws.Range("I21").Value2 = 1790
$ws.Range("J21").Value2 = 1891
$ws.Range("K21").Value2 = -5.341089370703
$ws.Range("L21").Value2 = -2.770233568712
$ws.Range("M21").Value2 = 24.133148404993
$ws.Range("N21").Value2 = -73.402674591381

# Row 22
$ws.Range("C22").Value2 = 3
$ws.Range("E22").Value2 = 0
$ws.Range("F22").Value2 = 7
$ws.Range("G22").Value2 = 14
$ws.Range("H22").Value2 = -50
$ws.Range("I22").Value2 = 95
$ws.Range("J22").Value2 = 115
$ws.Range("K22").Value2 = -17.391304347826
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = 35.714285714285

# Row 24
$ws.Range("C24").Value2 = 48
$ws.Range("D24").Value2 = 67
$ws.Range("E24").Value2 = -28.358208955223
$ws.Range("F24").Value2 = 265
$ws.Range("G24").Value2 = 327
$ws.Range("H24").Value2 = -18.960244648318
$ws.Range("I24").Value2 = 3850
$ws.Range("J24").Value2 = 4169
$ws.Range("K24").Value2 = -7.651715039577
$ws.Range("L24").Value2 = -4.561229548834
$ws.Range("M24").Value2 = 127.40696987596

# Row 25
$ws.Range("C25").Value2 = 48
$ws.Range("D25").Value2 = 77
$ws.Range("E25").Value2 = -37.662337662337
$ws.Range("F25").Value2 = 266
$ws.Range("G25").Value2 = 323
$ws.Range("H25").Value2 = -17.647058823529
$ws.Range("I25").Value2 = 3731
$ws.Range("J25").Value2 = 4182
$ws.Range("K25").Value2 = -10.784313725490
$ws.Range("L25").Value2 = -7.717041800643

# Row 26
$ws.Range("C26").Value2 = 6
$ws.Range("D26").Value2 = 7
$ws.Range("E26").Value2 = -14.285714285714
$ws.Range("F26").Value2 = 30
$ws.Range("G26").Value2 = 35
$ws.Range("H26").Value2 = -14.285714285714
$ws.Range("I26").Value2 = 447
$ws.Range("J26").Value2 = 413
$ws.Range("K26").Value2 = 8.232445520581
$ws.Range("L26").Value2 = 19.518716577540
$ws.Range("M26").Value2 = 64.944649446494

# Row 28
$ws.Range("G28").Value2 = 5
$ws.Range("H28").Value2 = -40
$ws.Range("I28").Value2 = 96
$ws.Range("J28").Value2 = 96
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = -15.789473684210

# Row 31
$ws.Range("G31").Value2 = 1

# --- Column width adjustments for columns I and J (9 and 10) ---
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
